# Update cryptos list (Price / Volume(1h) columns) per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D keeps its original text format (values like
# "26.671.41", "0.9983" or "242.00" must stay as text, not be
# reinterpreted/normalized as numbers by Excel).
$ws.Range("D2:D51").NumberFormat = "@"

$updates = @{
    2  = @{ D = "26.671.41";    E = "  +0.68%  " }
    3  = @{ D = "1.726.42";     E = "  -0.29%  " }
    4  = @{ D = "0.9983";       E = "  -0.28%  " }
    5  = @{ D = "242.00";       E = "  -0.88%  " }
    6  = @{ D = "0.9988";       E = "  -0.26%  " }
    7  = @{ D = "0.4930";       E = "  +0.48%  " }
    8  = @{ D = "0.2617";       E = "  -0.31%  " }
    9  = @{ D = "0.06226";      E = "  +0.68%  " }
    10 = @{ D = "1.731.13";     E = "  -0.15%  " }
    11 = @{ D = "15.86";        E = "  +2.47%  " }
    12 = @{ D = "0.06994";      E = "  -0.27%  " }
    13 = @{ D = "0.6109";       E = "  +1.60%  " }
    14 = @{ D = "4.508";        E = "  -0.85%  " }
    15 = @{ D = "77.29";        E = "  -0.23%  " }
    16 = @{ D = "0.9985";       E = "  -0.28%  " }
    17 = @{ D = "26.493.00";    E = "  -0.06%  " }
    18 = @{ D = "0.9985";       E = "  -0.29%  " }
    19 = @{ D = "0.000007214";  E = "  +1.86%  " }
    20 = @{ D = "11.42";        E = "  -0.01%  " }
    21 = @{ D = "1.953.07";     E = "  -0.28%  " }
    22 = @{ D = "4.466";        E = "  -0.26%  " }
    23 = @{ D = "8.566";        E = "  -0.42%  " }
    24 = @{ D = "5.103";        E = "  -1.49%  " }
    25 = @{ D = "138.02";       E = "  -0.54%  " }
    26 = @{ E = "  +0.95%  " }
    27 = @{ D = "1.766";        E = "  +2.97%  " }
    28 = @{ D = "1.382";        E = "  -2.33%  " }
    29 = @{ D = "106.23";       E = "  -0.28%  " }
    30 = @{ D = "3.918";        E = "  -1.53%  " }
    31 = @{ D = "0.07986";      E = "  +0.20%  " }
    32 = @{ D = "3.670";        E = "  -0.69%  " }
    33 = @{ D = "0.04508";      E = "  -0.25%  " }
    34 = @{ D = "0.9978";       E = "  -0.28%  " }
    35 = @{ D = "2.611";        E = "  -0.15%  " }
    36 = @{ D = "1.001";        E = "  +0.08%  " }
    37 = @{ D = "0.6264";       E = "  +0.24%  " }
    38 = @{ D = "0.9337";       E = "  +3.15%  " }
    39 = @{ D = "2.045";        E = "  +2.42%  " }
    40 = @{ D = "2.415";        E = "  +0.20%  " }
    41 = @{ D = "0.9998";       E = "  -0.27%  " }
    42 = @{ D = "0.01516";      E = "  +1.90%  " }
    43 = @{ D = "5.570";        E = "  +1.89%  " }
    44 = @{ D = "99.40";        E = "  -1.11%  " }
    45 = @{ D = "0.3858";       E = "  -0.15%  " }
    46 = @{ D = "6.895";        E = "  +3.85%  " }
    47 = @{ D = "0.1156";       E = "  +0.15%  " }
    48 = @{ D = "0.05381";      E = "  +0.28%  " }
    49 = @{ D = "7.893";        E = "  +2.87%  " }
    50 = @{ D = "30.29";        E = "  -0.05%  " }
    51 = @{ D = "51.74";        E = "  +1.37%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        $ws.Cells.Item($row, 4).Value = $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $ws.Cells.Item($row, 5).Value = $vals["E"]
    }
}
